$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "72.403.71"
$ws.Range("E2").Value = "  +4.30%  "

# Row 3
$ws.Range("D3").Value = "4.044.89"
$ws.Range("E3").Value = "  +3.57%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "518.29"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.04%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.21%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.721"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +17.65%  "

# Row 8
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.756"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +4.88%  "

# Row 10
$ws.Range("E10").Value = "  +1.35%  "

# Row 11
$ws.Range("E11").Value = "  -2.46%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.01"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +11.49%  "

# Row 13
$ws.Range("E13").Value = "  +6.27%  "

# Row 14
$ws.Range("D14").Value = "4.690.40"
$ws.Range("E14").Value = "  +3.46%  "

# Row 15
$ws.Range("D15").Value = "4.064.47"
$ws.Range("E15").Value = "  +4.33%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.06"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +6.64%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.08"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.11%  "

# Row 18
$ws.Range("E18").Value = "  -1.39%  "

# Row 19
$ws.Range("E19").Value = "  -1.88%  "

# Row 20
$ws.Range("D20").Value = "72.274.87"
$ws.Range("E20").Value = "  +4.18%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "441.24"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.69%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "103.89"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +17.21%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.60"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +6.16%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.58"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.60%  "

# Row 25
$ws.Range("E25").Value = "  -0.77%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.44"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.45%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.01"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.83%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.81"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.94%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.82"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.53%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.11"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +10.01%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.61"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.23%  "

# Row 32
$ws.Range("E32").Value = "  +2.22%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "677.87"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.26%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.82"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +13.69%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "67.38"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.89%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "42.42"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.71%  "

# Row 37
$ws.Range("B37").Value = "TheGraph"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.431"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.32%  "

# Row 38
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0862"
$ws.Range("E38").Value = "  +1.77%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.55"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +11.67%  "

# Row 40
$ws.Range("E40").Value = "  +0.55%  "

# Row 41
$ws.Range("E41").Value = "  +0.18%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0494"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.77%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.13%  "

# Row 44
$ws.Range("E44").Value = "  +2.95%  "

# Row 45
$ws.Range("E45").Value = "  +12.21%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.73"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.84%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.46"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.11%  "

# Row 48
$ws.Range("E48").Value = "  +2.29%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.03"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +6.44%  "

# Row 50
$ws.Range("E50").Value = "  +1.35%  "

# Row 51
$ws.Range("E51").Value = "  +1.15%  "
